# Update "1 Campaigns for Import.xlsx" - "My Campaigns" sheet
# The sample Campaign rows (D2:E18) are relabeled from coffee-equipment /
# generic marketing examples to a refreshed set of campaign names + types
# (including new bike-themed campaigns), the table is restyled, the Name
# column is narrowed, and the active selection is moved to the Campaign
# Type column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("My Campaigns")

# --- Row 2 ---------------------------------------------------------------
$ws.Range("D2").Value = "Customer Reference Lead"
$ws.Range("E2").Value = "Other"

# --- Row 3 ---------------------------------------------------------------
$ws.Range("D3").Value = "Market Trends Newsletter"
$ws.Range("E3").Value = "Direct Marketing"

# --- Row 4 ---------------------------------------------------------------
$ws.Range("D4").Value = "Monthly Newsletter"
$ws.Range("E4").Value = "Direct Marketing"

# --- Row 5 ---------------------------------------------------------------
$ws.Range("D5").Value = "New Product Releases"
$ws.Range("E5").Value = "Direct Marketing"

# --- Row 6 ---------------------------------------------------------------
$ws.Range("D6").Value = "Search Results"
$ws.Range("E6").Value = "Advertisement"

# --- Row 7 ---------------------------------------------------------------
$ws.Range("D7").Value = "Customer Care Campaign"
$ws.Range("E7").Value = "Other"

# --- Row 8 ---------------------------------------------------------------
$ws.Range("D8").Value = "Customer Follow-up"
$ws.Range("E8").Value = "Direct Marketing"

# --- Row 9 ---------------------------------------------------------------
$ws.Range("D9").Value = "Commercial Tradeshow"
$ws.Range("E9").Value = "Event"

# --- Row 10 --------------------------------------------------------------
$ws.Range("D10").Value = "Consumer Tradeshow"
$ws.Range("E10").Value = "Event"

# --- Row 11 --------------------------------------------------------------
$ws.Range("D11").Value = "In-App Video Placement"
$ws.Range("E11").Value = "Advertisement"

# --- Row 12 --------------------------------------------------------------
$ws.Range("D12").Value = "QuarterlySales Contest"
$ws.Range("E12").Value = "Other"

# --- Row 13 --------------------------------------------------------------
$ws.Range("D13").Value = "Le Tour Bundle"
$ws.Range("E13").Value = "Co-branding"

# --- Row 14 --------------------------------------------------------------
$ws.Range("D14").Value = "Family Bike Collection"
$ws.Range("E14").Value = "Co-branding"

# --- Row 15 --------------------------------------------------------------
$ws.Range("D15").Value = "Commuter Special"
$ws.Range("E15").Value = "Co-branding"

# --- Row 16 --------------------------------------------------------------
$ws.Range("D16").Value = "Tube Advert Special"
$ws.Range("E16").Value = "Co-branding"

# --- Row 17 --------------------------------------------------------------
$ws.Range("D17").Value = "ProAm Collection"
$ws.Range("E17").Value = "Co-branding"

# --- Row 18 --------------------------------------------------------------
$ws.Range("D18").Value = "Game Sponsorship "
$ws.Range("E18").Value = "Co-branding"

# --- Table style -----------------------------------------------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.TableStyle = "TableStyleMedium3"

# --- Column width (Name column narrowed from 62 to ~31.18 characters) -----
$ws.Columns.Item(4).ColumnWidth = 30.33

# --- Selection moved to the Campaign Type column ---------------------------
$ws.Range("E2:E18").Select()
